$d = $word.ActiveDocument

# --- Locate the "Pull" paragraph (unique, whole-word match) ---
$rng = $d.Content
$found = $rng.Find.Execute("Pull", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find the 'Pull' paragraph" }
$pullPara = $rng.Paragraphs(1)
$exPara = $pullPara.Next()

if ($exPara.Range.Text.Substring(0,3) -ne "Ex:") { throw "Unexpected paragraph after 'Pull': $($exPara.Range.Text)" }

# --- 1) Turn the old "Pull" paragraph into the new big heading ---
$headingXml = @'
<w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:color w:val="984806" w:themeColor="accent6" w:themeShade="80"/>
          <w:sz w:val="52"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:color w:val="984806" w:themeColor="accent6" w:themeShade="80"/>
          <w:sz w:val="52"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t xml:space="preserve">Best Complete example </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:color w:val="984806" w:themeColor="accent6" w:themeShade="80"/>
          <w:sz w:val="52"/>
        </w:rPr>
        <w:t>(Windows)</w:t>
      </w:r>
    </w:p>
'@
$pullPara.Range.InsertXML($headingXml)

# --- 2) Insert the whole new Dockerfile walkthrough block, plus a relocated
#        (bookmark-free, page-break-free) "Pull" paragraph, right before "Ex:" ---
$insertXml = @'
<w:p>
      <w:pPr>
        <w:rPr>
          <w:color w:val="FF0000"/>
          <w:sz w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:color w:val="FF0000"/>
          <w:sz w:val="28"/>
        </w:rPr>
        <w:t>Dockerfile</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:color w:val="FF0000"/>
          <w:sz w:val="28"/>
        </w:rPr>
        <w:t>:</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
        </w:rPr>
        <w:t xml:space="preserve">#simple example of </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
        </w:rPr>
        <w:t>dockerfile</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p><w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
        </w:rPr>
        <w:t xml:space="preserve">FROM </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
        </w:rPr>
        <w:t>ubuntu</w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
        </w:rPr>
        <w:t>:latest</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:proofErr w:type="gramEnd"/>
    </w:p><w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
        </w:rPr>
        <w:t>MAINT</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
        </w:rPr>
        <w:t>AINER singhsushil5277@gmail.com</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
        </w:rPr>
        <w:t>RUN apt-get update</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
        </w:rPr>
        <w:t>RUN apt-get install -y python</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="28"/>
        </w:rPr>
      </w:pPr>
    </w:p><w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
        </w:rPr>
        <w:t>COPY hello.py /home</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
        </w:rPr>
        <w:t>WORKDIR /home</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:rPr>
          <w:color w:val="FF0000"/>
          <w:sz w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:color w:val="FF0000"/>
          <w:sz w:val="28"/>
        </w:rPr>
        <w:t>Bulit</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:color w:val="FF0000"/>
          <w:sz w:val="28"/>
        </w:rPr>
        <w:t xml:space="preserve"> command:</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>docker</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> build </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> –</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>t</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> "</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>hello_python:dockerfiles</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve">" </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>C:\Users\admin\Desktop\pythondockerfiles</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:rPr>
          <w:color w:val="FF0000"/>
          <w:sz w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:color w:val="FF0000"/>
          <w:sz w:val="28"/>
        </w:rPr>
        <w:t>Run Command:</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
        </w:rPr>
        <w:t>docker</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
        </w:rPr>
        <w:t xml:space="preserve"> run -it </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
        </w:rPr>
        <w:t>hello_python:dockerfiles</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
        </w:rPr>
        <w:t xml:space="preserve"> python hello.py</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
    </w:p><w:p>
      <w:pPr>
        <w:rPr>
          <w:sz w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
        </w:rPr>
        <w:t>Pull</w:t>
      </w:r>
    </w:p>
'@
$insertPoint = $exPara.Range
$insertPoint.Collapse(1)
$insertPoint.InsertXML($insertXml)

# --- 3) The old bookmark ("_GoBack") used to sit at the end of the
#        "-> repo name" paragraph; it now lives in its own paragraph above,
#        so strip it from the "repo name" paragraph's content. ---
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("repo name", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) { throw "Could not find the 'repo name' paragraph" }
$repoPara = $rng2.Paragraphs(1)
$repoXml = @'
<w:p>
      <w:pPr>
        <w:rPr>
          <w:color w:val="FF0000"/>
          <w:sz w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
        </w:rPr>
        <w:tab/>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:color w:val="FF0000"/>
          <w:sz w:val="28"/>
        </w:rPr>
        <w:t>sk5277/</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:color w:val="FF0000"/>
          <w:sz w:val="28"/>
        </w:rPr>
        <w:t>fintech</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:color w:val="FF0000"/>
          <w:sz w:val="28"/>
        </w:rPr>
        <w:t>-&gt; repo name</w:t>
      </w:r>
    </w:p>
'@
$repoPara.Range.InsertXML($repoXml)

Write-Output "edit applied"
